# "Generate Report for Handback" — localization-status.xlsx
#
# a.md has come back from localization for zh-cn and de-de: the Overview /
# per-language Status cells flip from "Ready for handoff" to
# "Handed back: in sync with en-US", the per-language rows gain a
# "Latest Target File" / "Latest Handback File" / "Latest Handback DateTime"
# (a.md is now hyperlinked out of column I), and the Status/Handback-File
# columns widen to fit the new text.

$wb = $excel.ActiveWorkbook

$aUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ee29a162cb106bb16bd888368e03f645eb545fc0/e2e/a.md"

# ---------------------------------------------------------------------
# Overview sheet: zh-cn / de-de status columns (E, F) for both rows.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"
$wsOverview.Columns.Item(5).ColumnWidth = 29.1
$wsOverview.Columns.Item(6).ColumnWidth = 29.1

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value = "Handed back: in sync with en-US"
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Columns.Item(3).ColumnWidth = 29.1

# Latest Target File (I) now points at a.md, like column A.
$wsZh.Range("I2").Value = "a.md"
$wsZh.Range("I3").Value = "a.md"
$wsZh.Range("I2").Style = "HyperLink"
$wsZh.Range("I3").Style = "HyperLink"
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

# Latest Handback File (J) now has the generated xlf.
$wsZh.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.zh-cn.xlf"
$wsZh.Columns.Item(10).ColumnWidth = 39.2

# Latest Handback DateTime (K) — was the zero date, now the real handback time.
$wsZh.Range("K2").Value = "2016-08-26 06:35:42"
$wsZh.Range("K3").Value = "2016-08-26 06:35:42"

# ---------------------------------------------------------------------
# de-de sheet (same shape of edit, different xlf / timestamp)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Columns.Item(3).ColumnWidth = 29.1

$wsDe.Range("I2").Value = "a.md"
$wsDe.Range("I3").Value = "a.md"
$wsDe.Range("I2").Style = "HyperLink"
$wsDe.Range("I3").Style = "HyperLink"
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $aUrl, [Type]::Missing, [Type]::Missing, "a.md")

$wsDe.Range("J2").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Range("J3").Value = "a.6631f68b315a3f7ddcdc141802fdb6ef151d1df2.de-de.xlf"
$wsDe.Columns.Item(10).ColumnWidth = 39.2

$wsDe.Range("K2").Value = "2016-08-26 06:35:48"
$wsDe.Range("K3").Value = "2016-08-26 06:35:48"
